$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 44509
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 700
$ws.Range("Q6").Value = 25

# Row 7
$ws.Range("D7").Value = 44159
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("J7").Value = 28
$ws.Range("K7").Value = 19000
$ws.Range("L7").Value = 19500
$ws.Range("M7").Value = 19268
$ws.Range("N7").Value = "$/saco 25 kilos"
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 771
$ws.Range("Q7").Value = 25

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("J8").Value = 56
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19464
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 779
$ws.Range("Q8").Value = 25

# Row 9
$ws.Range("D9").Value = 44260
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("J9").Value = 33
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 23000
$ws.Range("M9").Value = 22545
$ws.Range("N9").Value = "$/saco 30 kilos"
$ws.Range("O9").Value = "Región de La Araucanía"
$ws.Range("P9").Value = 22545
$ws.Range("Q9").Value = 1

# Row 10
$ws.Range("D10").Value = 44524
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 500
$ws.Range("Q10").Value = 25

# Row 11
$ws.Range("D11").Value = 44216
$ws.Range("H11").Value = "Perfection"
$ws.Range("J11").Value = 43
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24419
$ws.Range("N11").Value = "$/saco 30 kilos"
$ws.Range("O11").Value = "Región de La Araucanía"
$ws.Range("P11").Value = 24419
$ws.Range("Q11").Value = 1

# Row 12
$ws.Range("D12").Value = 44186
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 21000
$ws.Range("L12").Value = 22000
$ws.Range("M12").Value = 21500
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Provincia de Diguillín"
$ws.Range("P12").Value = 860
$ws.Range("Q12").Value = 25

# Row 13
$ws.Range("D13").Value = 44215
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("J13").Value = 42
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19524
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Región de La Araucanía"
$ws.Range("P13").Value = 781
$ws.Range("Q13").Value = 25

# Row 14
$ws.Range("D14").Value = 44544
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Provincia de Diguillín"
$ws.Range("P14").Value = 540
$ws.Range("Q14").Value = 25

# Row 15
$ws.Range("D15").Value = 44530
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Diguillín"
$ws.Range("P15").Value = 580
$ws.Range("Q15").Value = 25

# Row 16
$ws.Range("D16").Value = 44530
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 12500
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 500
$ws.Range("Q16").Value = 25

# Row 17
$ws.Range("D17").Value = 44537
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 13500
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 540
$ws.Range("Q17").Value = 25

# Row 18
$ws.Range("D18").Value = 44546
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 540
$ws.Range("Q18").Value = 25

# Row 19
$ws.Range("D19").Value = 44162
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 18500
$ws.Range("L19").Value = 19000
$ws.Range("M19").Value = 18820
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 753
$ws.Range("Q19").Value = 25

# Row 20
$ws.Range("D20").Value = 44166
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17479
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 699
$ws.Range("Q20").Value = 25

# Row 21
$ws.Range("D21").Value = 44174
$ws.Range("H21").Value = "Perfection"
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19500
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Región de Ñuble"
$ws.Range("P21").Value = 780
$ws.Range("Q21").Value = 25

# Row 22
$ws.Range("D22").Value = 44273
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 22000
$ws.Range("M22").Value = 21091
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región de La Araucanía"
$ws.Range("P22").Value = 844
$ws.Range("Q22").Value = 25

# Row 23
$ws.Range("D23").Value = 44525
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12500
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 500
$ws.Range("Q23").Value = 25

# Row 24
$ws.Range("D24").Value = 44165
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("J24").Value = 42
$ws.Range("K24").Value = 18000
$ws.Range("L24").Value = 19000
$ws.Range("M24").Value = 18595
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 744
$ws.Range("Q24").Value = 25

# Row 25
$ws.Range("D25").Value = 44267
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = 24000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 24581
$ws.Range("N25").Value = "$/saco 30 kilos"
$ws.Range("O25").Value = "Región de La Araucanía"
$ws.Range("P25").Value = 24581
$ws.Range("Q25").Value = 1

# Row 26
$ws.Range("D26").Value = 44516
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 15500
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 620
$ws.Range("Q26").Value = 25

# Row 27
$ws.Range("D27").Value = 44176
$ws.Range("H27").Value = "Perfection"
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 21000
$ws.Range("M27").Value = 20500
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Provincia de Diguillín"
$ws.Range("P27").Value = 820
$ws.Range("Q27").Value = 25

# Row 28
$ws.Range("D28").Value = 44532
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 13000
$ws.Range("L28").Value = 14000
$ws.Range("M28").Value = 13500
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 540
$ws.Range("Q28").Value = 25

# Row 29
$ws.Range("D29").Value = 44498
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 17500
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 700
$ws.Range("Q29").Value = 25

# Row 30
$ws.Range("D30").Value = 44167
$ws.Range("H30").Value = "Perfection"
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 18000
$ws.Range("L30").Value = 19000
$ws.Range("M30").Value = 18500
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 740
$ws.Range("Q30").Value = 25

# Row 31
$ws.Range("D31").Value = 44161
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("J31").Value = 33
$ws.Range("K31").Value = 19000
$ws.Range("L31").Value = 19500
$ws.Range("M31").Value = 19303
$ws.Range("N31").Value = "$/saco 25 kilos"
$ws.Range("O31").Value = "Provincia de Diguillín"
$ws.Range("P31").Value = 772
$ws.Range("Q31").Value = 25

# Row 32
$ws.Range("D32").Value = 44161
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("J32").Value = 34
$ws.Range("K32").Value = 19500
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 19735
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 789
$ws.Range("Q32").Value = 25

# Row 33
$ws.Range("D33").Value = 44160
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("J33").Value = 62
$ws.Range("K33").Value = 19000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 19516
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 781
$ws.Range("Q33").Value = 25

# Row 34
$ws.Range("D34").Value = 44551
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15500
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Provincia de Diguillín"
$ws.Range("P34").Value = 620
$ws.Range("Q34").Value = 25
